# "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta"
#
# The EC (Estado de Cuenta) sheet currently lists two overdue-period rows
# for worker GLORIA PAOLA MARTIN BUELVAS (periods 2507 and 2506) along
# with totals for "VALOR MORA" and "Cant. Periodos" that reflect both
# rows. The update removes the older 2506 period row (keeping only the
# still-relevant 2507 row) and refreshes the two summary figures to
# match the now-single-period data set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# VALOR MORA total (E11): was the sum for both periods (88000), now just
# the remaining 2507 period's value.
$ws.Range("E11").Value = 80000

# Cant. Periodos (F13): was 2 overdue periods, now only 1.
$ws.Range("F13").Value = 1

# Drop the 2506 period detail row (row 17: CC 1143328792, GLORIA PAOLA
# MARTIN BUELVAS, period 2506, 8000, 2000000). Deleting the whole row
# shifts the trailing signature block (previously rows 22-23) up to
# rows 21-22.
$ws.Rows("17").Delete()
